$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 555; this pushes the existing rows 555..675
# down to 556..676, matching the target diff (which shows the original
# row 555 data re-appearing at row 556, and so on, with the original
# row 675 re-appearing at the new row 676).
$ws.Rows.Item(555).Insert()

# Populate the newly inserted row 555 with its data. Columns A, B, C,
# E, F, G, H, I, N, Q, R follow the same static pattern used by every
# other data row in this sheet.
$ws.Cells.Item(555, 1).Value = 4
$ws.Cells.Item(555, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(555, 3).Value = "Los Lagos"

# Column D holds a date serial; copy the date number format used by
# the surrounding rows (style index 2 / numFmtId 165) from row 556.
$ws.Cells.Item(555, 4).Value = 45244
$ws.Cells.Item(555, 4).NumberFormat = $ws.Cells.Item(556, 4).NumberFormat

$ws.Cells.Item(555, 5).Value = 10
$ws.Cells.Item(555, 6).Value = 100114013
$ws.Cells.Item(555, 7).Value = "Zanahoria"
$ws.Cells.Item(555, 8).Value = "Sin especificar"
$ws.Cells.Item(555, 9).Value = "Primera"
$ws.Cells.Item(555, 10).Value = 700
$ws.Cells.Item(555, 11).Value = 9000
$ws.Cells.Item(555, 12).Value = 10000
$ws.Cells.Item(555, 13).Value = 9500
$ws.Cells.Item(555, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(555, 15).Value = "Región Metropolitana"
$ws.Cells.Item(555, 16).Value = 475
$ws.Cells.Item(555, 17).Value = 20
$ws.Cells.Item(555, 18).Value = "Hortaliza"
